$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# NALCO published a new circular (22-11-2025). Insert it as the new top
# data row (row 2) and push the existing 18 data rows (2..19) down to
# (3..20), keeping all of their original values intact.
# ---------------------------------------------------------------------

# Stash copies of the plain data-row cell styles used by this table (no
# borders, centered, "General"/0.000 number format) in scratch cells far
# outside the printed table so they can be re-applied later without
# fear of being clobbered by in-between edits.
$ws.Range("E2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats -> plain text-ish style (s=3)
$ws.Range("D2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # xlPasteFormats -> price style (s=4)
$excel.CutCopyMode = $false

# 1) Give the brand-new row 20 the same look (borders/number format) as
#    the rest of the table before any values land there.
$ws.Range("A19:F19").Copy()
$ws.Range("A20:F20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Shift the 18 existing data rows down by one, bottom-up so a row is
#    never overwritten before it has been read.
$slNo  = @(18,17,16,15,14,13,12,11,10,9,8,7,6,5,4,3,2,1)
$price = @(281.95,292.65,297.15,294.05,288.55,282.45,285.05,282.85,277.95,274.95,270.25,275.25,278.95,272.05,271.05,264.35,269.45,268.25)
$date  = @("19-11-2025","07-11-2025","01-11-2025","30-10-2025","25-10-2025","17-10-2025","14-10-2025","09-10-2025","01-10-2025","30-09-2025","25-09-2025","20-09-2025","17-09-2025","01-09-2025","28-08-2025","21-08-2025","15-08-2025","07-08-2025")
$link  = @(
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-19-11-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
)

# Column E holds circular dates as plain TEXT (not Excel date serials),
# matching the source file. Force text format first so ambiguous-looking
# strings like "07-11-2025" aren't auto-converted to dates on entry.
$ws.Range("E2:E20").NumberFormat = "@"

for ($i = 17; $i -ge 0; $i--) {
    $destRow = $i + 3   # old row (i+2) -> new row (i+3)
    $ws.Cells.Item($destRow, 1).Value2 = $slNo[$i]
    $ws.Cells.Item($destRow, 2).Value2 = "ALUMINIUM INGOT"
    $ws.Cells.Item($destRow, 3).Value2 = "IE07"
    $ws.Cells.Item($destRow, 4).Value2 = $price[$i]
    $ws.Cells.Item($destRow, 5).Value2 = $date[$i]
    $ws.Cells.Item($destRow, 6).Value2 = $link[$i]
}

# 3) New row 2: the latest circular.
$ws.Cells.Item(2, 1).Value2 = 19
$ws.Cells.Item(2, 2).Value2 = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value2 = "IE07"
$ws.Cells.Item(2, 4).Value2 = 283.55
$ws.Cells.Item(2, 5).Value2 = "22-11-2025"
$newLink = "https://nalcoindia.com/wp-content/uploads/2025/11/Ingot-22-11-2025.pdf"
$ws.Cells.Item(2, 6).Value2 = $newLink

# Restore the plain table look on columns D/E that the steps above may
# have disturbed (text-number-format stamp on E, float noise is fine).
$ws.Range("Z1").Copy()
$ws.Range("E2:E20").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("D2:D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4) Rebuild the hyperlinks on column F for rows 2..20 top to bottom, so
#    the relationship ids come out sequentially as rId1..rId19. (Deleting
#    any single cell's Hyperlinks collection clears the whole sheet's, so
#    this one call is enough to start from a clean slate.)
$ws.Cells.Item(2, 6).Hyperlinks.Delete()

$allLinks = @($newLink) + $link
for ($r = 2; $r -le 20; $r++) {
    $u = $allLinks[$r - 2]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $u)
}

# Adding a hyperlink re-styles the touched cell with the blue/underline
# "Hyperlink" look; this sheet never used that look for its links, so
# stamp the plain data-row format back onto all of column F's data cells.
$ws.Range("Z1").Copy()
$ws.Range("F2:F20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Clean up scratch cells so they don't show up in the saved sheet.
$ws.Range("Z1:Z2").Clear()

$wb.Save()
